# Contoso Chai Tea market trends 2023 - row 6 data correction.
#
# D6 (Ventes de chaï préconfectionné (unités)) and E6 (Engagement sur les
# réseaux sociaux (vues)) are changed from numeric values (436 / 1705) to
# text time-stamps ("4:36" / "05:17"). Column B holds a shared formula
# (=SUM(C+D)) that would error out (#VALUE!) if Excel recalculated it
# against the new text value in D6, but the source edit only touched the
# raw cell contents, so B6's cached result must stay exactly as it was.
# Switching to manual calculation before writing the new values (and
# leaving it there) prevents that ripple/recalculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$excel.Calculation = -4135   # xlCalculationManual - avoid recalculating B6

$ws.Range("D6").Value = "4:36"
$ws.Range("E6").Value = "05:17"
